$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before B; this shifts the old B,C,D -> C,D,E
$ws.Columns("B:B").Insert()

# New header row: C1 (caso_general) must be written before B1 (tipo) so that
# the shared-string table picks up "caso_general" at the lower index, matching
# the original author's save order.
$ws.Range("C1").Value = "caso_general"
$ws.Range("B1").Value = "tipo"

# New "tipo" column values classifying each variable group
$ws.Range("B9").Value = "t_sexo"
$ws.Range("B10").Value = "t_edad"
$ws.Range("B11:B22").Value = "t_alzanzado"

# Style B1: bold (inherited), centered, with a yellow fill
$ws.Range("B1").Interior.Color = 65535
$ws.Range("B1").HorizontalAlignment = -4108

# Column widths: B matches A's (bestFit) width as closely as this host allows;
# C:E get a fixed 21.5-character width.
$ws.Columns("B:B").ColumnWidth = 9.5
$ws.Columns("C:E").ColumnWidth = 20.6666666666667

# Move the active selection to C13 (matches the saved selection in the workbook)
[void]$ws.Range("C13").Select()
